$d = $word.ActiveDocument
$rng = $d.Content

function Do-Replace($findText, $replaceText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $r.Find.Replacement.ClearFormatting()
    $ok = $r.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) { throw "Replace failed for: $findText" }
}

# Phase 1: move original texts to unique placeholders
Do-Replace ('Apresentar conceitos e metodologias fundamentais para concepção e projeto de fábricas, com ênfase no projeto e organização dos processos de produção, movimentação e armazenagem de materiais, visando o adequado aproveitamento do espaço físico e a eficiência nos fluxos de materiais e ordens de produção no ambiente interno das fábricas.') '@@SWAP_9@@'
Do-Replace ('To introduce fundamental concepts and methodologies for plant design, with emphasis on the design and organization of production processes, materials handling and warehousing, aiming at appropriate use of space and efficient flow of materials and jobs in the internal environment of manufacturing plants') '@@SWAP_10@@'
Do-Replace ('8971158 - Claudemir Leif Tramarico') '@@SWAP_12@@'
Do-Replace ('Arranjo Físico da Fábrica; Sistema de Movimentação e Armazenagem de Materiais; Análise do Fluxo de Materiais. Planejamento do Layout') '@@SWAP_14@@'
Do-Replace ('Plant Layout; Materials Handling and Warehousing System; Materials Flow Analysis. Layout Planning') '@@SWAP_15@@'
Do-Replace ('1. Objetivos de Desempenho de Empresas de Manufatura e de suas Fábricas. 2. Conceitos de Produto, Recurso e Processo para Projeto da Fábrica. 3. Tipos de Produção e Tipos de Arranjo Físico. 4. Planejamento do Arranjo Físico e dos Fluxos Internos. 5. Manufatura Celular. 6. Planejamento do Sistema de Movimentação e Armazenagem de Materiais.7. Planejamento do Layout: Sistema SLP') '@@SWAP_17@@'
Do-Replace ('Provas, atividades em grupo e atividades individuais.') '@@SWAP_21@@'
Do-Replace ('Média das atividades avaliativas') '@@SWAP_23@@'
Do-Replace ('MF = (0,5 M + 0,5 R) M = Média de aproveitamento do aluno, antes da recuperação R = Nota de uma prova de recuperação MF = nota final de aproveitamento, após a recuperação Aprovação com média final de aproveitamento maior ou igual a 5,0. A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre. Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.') '@@SWAP_25@@'
Do-Replace ('BANZATO, Eduardo et al. Atualidades na armazenagem. São Paulo: IMAM, 2003.' + [char]11 + 'BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. São Paulo, Edgar Blücher, 1977.' + [char]11 + 'GURGEL, F.A.C. Administração de recursos materiais e patrimoniais. 2a. Edição. São Paulo. Editora Cengage. 2013. ' + [char]11 + 'FRANCISCHINI, P.G.; VALLE, C.E. Implantação de Indústrias. Rio de Janeiro, LTC Editora, 1975.' + [char]11 + 'LEE, Q et al. Projeto de Instalações e Locais de Trabalho. São Paulo: IMAM, 1998.' + [char]11 + 'MOURA, Reinaldo Aparecido. Sistemas e técnicas de movimentação e armazenagem de materiais. IMAM, 2012.' + [char]11 + 'NEWMANN, C.; SCALICE, R.K. Projeto de Fábrica e Layout. Rio de Janeiro, Elsevier, 2015.' + [char]11 + 'Müther, R. Planejamento do Layout: Sistema SLP. São Paulo, Edgard Blücher, 1978. ' + [char]11 + 'SLACK, Nigel et al. Administração da produção. São Paulo: Atlas, 8ª ed. 2018.' + [char]11 + 'TOMPKINS, James A. et al. Planejamento de instalações. Editora LTC:, 2013.') '@@SWAP_27@@'

# Phase 2: placeholders -> final texts
Do-Replace '@@SWAP_9@@' ('Arranjo Físico da Fábrica; Sistema de Movimentação e Armazenagem de Materiais; Análise do Fluxo de Materiais. Planejamento do Layout')
Do-Replace '@@SWAP_10@@' ('Plant Layout; Materials Handling and Warehousing System; Materials Flow Analysis. Layout Planning')
Do-Replace '@@SWAP_12@@' ('Apresentar conceitos e metodologias fundamentais para concepção e projeto de fábricas, com ênfase no projeto e organização dos processos de produção, movimentação e armazenagem de materiais, visando o adequado aproveitamento do espaço físico e a eficiência nos fluxos de materiais e ordens de produção no ambiente interno das fábricas.')
Do-Replace '@@SWAP_14@@' ('1. Objetivos de Desempenho de Empresas de Manufatura e de suas Fábricas. 2. Conceitos de Produto, Recurso e Processo para Projeto da Fábrica. 3. Tipos de Produção e Tipos de Arranjo Físico. 4. Planejamento do Arranjo Físico e dos Fluxos Internos. 5. Manufatura Celular. 6. Planejamento do Sistema de Movimentação e Armazenagem de Materiais.7. Planejamento do Layout: Sistema SLP')
Do-Replace '@@SWAP_15@@' ('To introduce fundamental concepts and methodologies for plant design, with emphasis on the design and organization of production processes, materials handling and warehousing, aiming at appropriate use of space and efficient flow of materials and jobs in the internal environment of manufacturing plants')
Do-Replace '@@SWAP_17@@' ('Provas, atividades em grupo e atividades individuais.')
Do-Replace '@@SWAP_21@@' ('Média das atividades avaliativas')
Do-Replace '@@SWAP_23@@' ('MF = (0,5 M + 0,5 R) M = Média de aproveitamento do aluno, antes da recuperação R = Nota de uma prova de recuperação MF = nota final de aproveitamento, após a recuperação Aprovação com média final de aproveitamento maior ou igual a 5,0. A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre. Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.')
Do-Replace '@@SWAP_25@@' ('BANZATO, Eduardo et al. Atualidades na armazenagem. São Paulo: IMAM, 2003.' + [char]11 + 'BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. São Paulo, Edgar Blücher, 1977.' + [char]11 + 'GURGEL, F.A.C. Administração de recursos materiais e patrimoniais. 2a. Edição. São Paulo. Editora Cengage. 2013. ' + [char]11 + 'FRANCISCHINI, P.G.; VALLE, C.E. Implantação de Indústrias. Rio de Janeiro, LTC Editora, 1975.' + [char]11 + 'LEE, Q et al. Projeto de Instalações e Locais de Trabalho. São Paulo: IMAM, 1998.' + [char]11 + 'MOURA, Reinaldo Aparecido. Sistemas e técnicas de movimentação e armazenagem de materiais. IMAM, 2012.' + [char]11 + 'NEWMANN, C.; SCALICE, R.K. Projeto de Fábrica e Layout. Rio de Janeiro, Elsevier, 2015.' + [char]11 + 'Müther, R. Planejamento do Layout: Sistema SLP. São Paulo, Edgard Blücher, 1978. ' + [char]11 + 'SLACK, Nigel et al. Administração da produção. São Paulo: Atlas, 8ª ed. 2018.' + [char]11 + 'TOMPKINS, James A. et al. Planejamento de instalações. Editora LTC:, 2013.')
Do-Replace '@@SWAP_27@@' ('8971158 - Claudemir Leif Tramarico')